# Updates the cryptos price/volume snapshot (and a few swapped coin rows)
# to match the latest scrape. Numeric-looking price strings are prefixed
# with a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr/text storage) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.928.60'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '3.740.02'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''601.48'
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = '''165.92'
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("D7").Value = '3.737.57'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("E10").Value = '  +4.18%  '
$ws.Range("D11").Value = '''6.42'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '''0.460'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '''37.80'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = '''0.0000248'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '4.366.50'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '3.745.15'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '69.015.07'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '''7.44'
$ws.Range("E18").Value = '  +2.26%  '
$ws.Range("D19").Value = '''17.72'
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").Value = '''0.114'
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").Value = '''11.14'
$ws.Range("E21").Value = '  +5.01%  '
$ws.Range("D22").Value = '''491.77'
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '''84.61'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").Value = '''0.0000148'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").Value = '''2.27'
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("D27").Value = '''12.35'
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = '''10.06'
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''8.23'
$ws.Range("E30").Value = '  +3.72%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").Value = '''2.44'
$ws.Range("E32").Value = '  -5.68%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''31.64'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("B34").Value = 'WrappedeETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D34").Value = '3.885.40'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").Value = '3.675.57'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '''5.93'
$ws.Range("E37").Value = '  +1.60%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("E39").Value = '  +5.51%  '
$ws.Range("D40").Value = '''3.19'
$ws.Range("E40").Value = '  +10.78%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '''0.325'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''2.00'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''48.58'
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").Value = '''426.16'
$ws.Range("E45").Value = '  -3.10%  '
$ws.Range("D46").Value = '''8.45'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D48").Value = '''40.04'
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("D49").Value = '''141.10'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").Value = '''1.30'
$ws.Range("E50").Value = '  +7.06%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.780.89'
$ws.Range("E51").Value = '  +0.17%  '
